# Gap-analysis sheet: add a second subject (DSA / 18CS32), recompute the
# SUM / Total / Gap% / Gaps-in-T summary rows for the now-two-subject table,
# un-merge the old A3:C4 placeholder block, widen column B, and tidy up the
# cell borders/alignment so every column (incl. the PSO columns that used to
# be blank) is boxed and centred the same way as the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- widen the "Subject Name" column (30 -> 40 chars) ----------------------
# ColumnWidth round-trips with a ~0.8333 char padding in this engine, so back
# that out to land exactly on width="40" in the saved xf.
$ws.Columns.Item(2).ColumnWidth = 39.166666666666664

# --- remove the old A3:C4 merge (row 3 becomes real per-subject data) ------
$ws.Range("A3:C4").UnMerge()

# --- row 3: second subject, "DSA" / "18CS32" --------------------------------
$ws.Range("B3").Value = "DSA"
$ws.Range("C3").Value = "18CS32"

function Set-RowValues($row, $startCol, $values) {
    $col = $startCol
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}

# --- row 4: SUM = A (was the old row 5 content) -----------------------------
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = "SUM = A"
$ws.Range("C4").Value = ""
Set-RowValues 4 4 @(5, 3.2, 3.2, 1, 2, 0, 0, 0, 4, 0, 0, 0, 4, 4)

# --- row 5: (Total no. of courses addressing each PO)= T -------------------
$ws.Range("B5").Value = "(Total no. of courses addressing each PO)= T"
Set-RowValues 5 4 @(2, 2, 2, 1, 1, 0, 0, 0, 2, 0, 0, 0, 2, 2)

# --- row 6: GAP G= (27-A)/(27))*100 -----------------------------------------
$ws.Range("A6").Value = ""
$ws.Range("B6").Value = "GAP G= (27-A)/(27))*100"
$ws.Range("C6").Value = ""
Set-RowValues 6 4 @(81.48148148148148, 88.14814814814815, 88.14814814814815, 96.29629629629629, 92.5925925925926, 100, 100, 100, 85.18518518518519, 100, 100, 100, 85.18518518518519, 85.18518518518519)

# --- row 7: Gaps in T --------------------------------------------------------
$ws.Range("B7").Value = "Gaps in T"
Set-RowValues 7 4 @(12, 12, 12, 13, 13, 14, 14, 14, 12, 14, 14, 14, 12, 12)

# --- row 8: wipe the old (narrower, 7-column) summary leftovers ------------
# give the whole row real cells first (H8:Q8 don't exist yet) so the
# ClearContents below leaves behind empty-but-styled cells across A:Q
$ws.Range("H8:Q8").Value = 0
$ws.Range("B2").Copy()
$ws.Range("A8:Q8").PasteSpecial(-4122)
$ws.Range("A8:Q8").ClearContents()

# --- style clean-up ----------------------------------------------------------
# B3 goes from "bordered + centred" to "bordered only" (matches B2's look)
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)

# D4:Q4, P6:Q6 and L7:Q7 pick up the bordered+centred look used everywhere
# else in the table (copy from A1, which already carries that style)
$ws.Range("A1").Copy()
$ws.Range("D4:Q4").PasteSpecial(-4122)
$ws.Range("P6:Q6").PasteSpecial(-4122)
$ws.Range("L7:Q7").PasteSpecial(-4122)

$excel.CutCopyMode = $false
